$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws 'D2' '42.085.10'
$ws.Range('E2').Value = '  -4.24%  '

# Row 3
Set-TextValue $ws 'D3' '2.229.88'
$ws.Range('E3').Value = '  -5.20%  '

# Row 4
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
Set-TextValue $ws 'D5' '242.73'
$ws.Range('E5').Value = '  +0.87%  '

# Row 6
$ws.Range('E6').Value = '  -6.15%  '

# Row 7
Set-TextValue $ws 'D7' '68.06'
$ws.Range('E7').Value = '  -7.72%  '

# Row 8
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
Set-TextValue $ws 'D9' '0.544'
$ws.Range('E9').Value = '  -9.25%  '

# Row 10
$ws.Range('E10').Value = '  -4.91%  '

# Row 11
Set-TextValue $ws 'D11' '58.26'
$ws.Range('E11').Value = '  -3.38%  '

# Row 12
Set-TextValue $ws 'D12' '35.46'
$ws.Range('E12').Value = '  +6.90%  '

# Row 13
$ws.Range('E13').Value = '  -2.97%  '

# Row 14
Set-TextValue $ws 'D14' '6.70'
$ws.Range('E14').Value = '  -8.04%  '

# Row 15
Set-TextValue $ws 'D15' '2.560.91'
$ws.Range('E15').Value = '  -5.08%  '

# Row 16
Set-TextValue $ws 'D16' '14.79'
$ws.Range('E16').Value = '  -8.36%  '

# Row 17
Set-TextValue $ws 'D17' '0.851'
$ws.Range('E17').Value = '  -6.22%  '

# Row 18
Set-TextValue $ws 'D18' '2.219.02'
$ws.Range('E18').Value = '  -5.46%  '

# Row 19
Set-TextValue $ws 'D19' '41.940.68'
$ws.Range('E19').Value = '  -4.30%  '

# Row 20
Set-TextValue $ws 'D20' '0.0₃0954'
$ws.Range('E20').Value = '  -7.80%  '

# Row 21
Set-TextValue $ws 'D21' '72.59'
$ws.Range('E21').Value = '  -7.61%  '

# Row 22
Set-TextValue $ws 'D22' '6.12'
$ws.Range('E22').Value = '  -8.12%  '

# Row 23
Set-TextValue $ws 'D23' '234.25'
$ws.Range('E23').Value = '  -7.36%  '

# Row 24
Set-TextValue $ws 'D24' '2.05'
$ws.Range('E24').Value = '  +11.59%  '

# Row 25
Set-TextValue $ws 'D25' '0.999'
$ws.Range('E25').Value = '  -0.20%  '

# Row 26
$ws.Range('E26').Value = '  -5.36%  '

# Row 27
$ws.Range('E27').Value = '  -2.87%  '

# Row 28
Set-TextValue $ws 'D28' '2.24'
$ws.Range('E28').Value = '  -3.05%  '

# Row 29
Set-TextValue $ws 'D29' '9.87'
$ws.Range('E29').Value = '  -5.69%  '

# Row 30
Set-TextValue $ws 'D30' '171.63'
$ws.Range('E30').Value = '  -2.45%  '

# Row 31
Set-TextValue $ws 'D31' '20.37'
$ws.Range('E31').Value = '  -8.65%  '

# Row 32
$ws.Range('E32').Value = '  -4.95%  '

# Row 33
$ws.Range('E33').Value = '  -6.36%  '

# Row 34
Set-TextValue $ws 'D34' '0.0712'
$ws.Range('E34').Value = '  -4.92%  '

# Row 35
Set-TextValue $ws 'D35' '5.19'
$ws.Range('E35').Value = '  -2.97%  '

# Row 36
Set-TextValue $ws 'D36' '4.66'
$ws.Range('E36').Value = '  -8.30%  '

# Row 37
Set-TextValue $ws 'D37' '3.88'
$ws.Range('E37').Value = '  +1.21%  '

# Row 38
Set-TextValue $ws 'D38' '23.01'
$ws.Range('E38').Value = '  +22.01%  '

# Row 39
Set-TextValue $ws 'D39' '0.0282'
$ws.Range('E39').Value = '  +3.82%  '

# Row 40
Set-TextValue $ws 'D40' '2.29'
$ws.Range('E40').Value = '  -4.06%  '

# Row 41
Set-TextValue $ws 'D41' '66.76'
$ws.Range('E41').Value = '  +2.87%  '

# Row 42
$ws.Range('E42').Value = '  -9.30%  '

# Row 43
Set-TextValue $ws 'D43' '9.03'
$ws.Range('E43').Value = '  -1.79%  '

# Row 44
Set-TextValue $ws 'D44' '4.91'
$ws.Range('E44').Value = '  -12.25%  '

# Row 45
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D45' '0.191'
$ws.Range('E45').Value = '  -4.11%  '

# Row 46
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D46' '0.100'
$ws.Range('E46').Value = '  -4.18%  '

# Row 47
$ws.Range('B47').Value = 'SynthetixNetwork'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue $ws 'D47' '4.63'
$ws.Range('E47').Value = '  +8.46%  '

# Row 48
$ws.Range('B48').Value = 'BinanceUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws 'D48' '1.00'
$ws.Range('E48').Value = '  -0.09%  '

# Row 49
Set-TextValue $ws 'D49' '1.19'
$ws.Range('E49').Value = '  -3.22%  '

# Row 50
Set-TextValue $ws 'D50' '2.81'
$ws.Range('E50').Value = '  -2.38%  '

# Row 51
$ws.Range('E51').Value = '  -4.63%  '
